$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fish species names to append below the existing header (A1)
$names = @(
    "Leptoscarus vaigiensis",
    "Lethrinus nebulosus",
    "Scarus ghobban",
    "Siganus canaliculutus",
    "Siganus sutor"
)

$row = 2
foreach ($name in $names) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $name
    $cell.Font.Name = "Monaco"
    $cell.Font.Size = 12
    $cell.Font.Color = 0
    $row = $row + 1
}

# Widen column A to fit the new content
$ws.Columns.Item(1).ColumnWidth = 30.33203125

# Update the active selection
$ws.Range("I7").Select()
